$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$addr, [string]$val) {
    $c = $ws.Range($addr)
    $escaped = $val -replace '"', '""'
    $c.Formula = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue 'D2' '51.701.71'
Set-TextValue 'E2' '  +1.27%  '
Set-TextValue 'D3' '3.028.17'
Set-TextValue 'E3' '  +2.33%  '
Set-TextValue 'E4' '  +0.03%  '
Set-TextValue 'D5' '380.53'
Set-TextValue 'E5' '  +0.25%  '
Set-TextValue 'D6' '102.85'
Set-TextValue 'E6' '  +0.56%  '
Set-TextValue 'D7' '0.547'
Set-TextValue 'E7' '  +0.84%  '
Set-TextValue 'D9' '0.594'
Set-TextValue 'E9' '  +1.50%  '
Set-TextValue 'D10' '36.84'
Set-TextValue 'E10' '  +1.39%  '
Set-TextValue 'E11' '  -0.34%  '
Set-TextValue 'E12' '  +1.11%  '
Set-TextValue 'D13' '3.509.32'
Set-TextValue 'E13' '  +2.41%  '
Set-TextValue 'D14' '18.53'
Set-TextValue 'E14' '  +0.58%  '
Set-TextValue 'D15' '7.74'
Set-TextValue 'E15' '  -0.83%  '
Set-TextValue 'D16' '3.025.00'
Set-TextValue 'E16' '  +2.69%  '
Set-TextValue 'E17' '  -4.25%  '
Set-TextValue 'E18' '  -15.46%  '
Set-TextValue 'D19' '51.706.57'
Set-TextValue 'E19' '  +1.29%  '
Set-TextValue 'D20' '3.09'
Set-TextValue 'E20' '  +0.27%  '
Set-TextValue 'E21' '  +1.02%  '
Set-TextValue 'D22' '0.0₃0962'
Set-TextValue 'E22' '  +1.06%  '
Set-TextValue 'E23' '  +0.66%  '
Set-TextValue 'D24' '268.35'
Set-TextValue 'E24' '  +0.71%  '
Set-TextValue 'E25' '  -5.98%  '
Set-TextValue 'D26' '8.26'
Set-TextValue 'E26' '  +2.61%  '
Set-TextValue 'D27' '7.62'
Set-TextValue 'E27' '  +8.81%  '
Set-TextValue 'E28' '  +4.91%  '
Set-TextValue 'E29' '  -0.10%  '
Set-TextValue 'D30' '26.23'
Set-TextValue 'E30' '  +1.79%  '
Set-TextValue 'D31' '0.109'
Set-TextValue 'E31' '  +0.88%  '
Set-TextValue 'D32' '10.29'
Set-TextValue 'E32' '  -0.31%  '
Set-TextValue 'B33' 'InjectiveProtocol'
Set-TextValue 'C33' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D33' '33.99'
Set-TextValue 'E33' '  -0.16%  '
Set-TextValue 'B34' 'OKB'
Set-TextValue 'C34' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D34' '50.50'
Set-TextValue 'E34' '  -0.08%  '
Set-TextValue 'B35' 'Toncoin'
Set-TextValue 'C35' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D35' '2.05'
Set-TextValue 'E35' '  -0.18%  '
Set-TextValue 'E36' '  +3.39%  '
Set-TextValue 'E37' '  -0.08%  '
Set-TextValue 'D38' '3.32'
Set-TextValue 'E38' '  +4.81%  '
Set-TextValue 'D39' '0.298'
Set-TextValue 'E39' '  +15.14%  '
Set-TextValue 'D40' '17.05'
Set-TextValue 'E40' '  +2.72%  '
Set-TextValue 'D41' '1.86'
Set-TextValue 'E41' '  +2.46%  '
Set-TextValue 'E42' '  +2.67%  '
Set-TextValue 'B43' 'Stellar'
Set-TextValue 'C43' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D43' '0.116'
Set-TextValue 'E43' '  -0.43%  '
Set-TextValue 'B44' 'NEARProtocol'
Set-TextValue 'C44' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D44' '3.79'
Set-TextValue 'E44' '  +5.90%  '
Set-TextValue 'B45' 'Monero'
Set-TextValue 'C45' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D45' '124.17'
Set-TextValue 'E45' '  +3.39%  '
Set-TextValue 'D46' '21.70'
Set-TextValue 'E46' '  +1.49%  '
Set-TextValue 'D47' '2.09'
Set-TextValue 'E47' '  +3.52%  '
Set-TextValue 'E48' '  +4.90%  '
Set-TextValue 'D49' '2.030.96'
Set-TextValue 'E49' '  +0.84%  '
Set-TextValue 'D50' '3.329.91'
Set-TextValue 'E50' '  +2.38%  '
Set-TextValue 'D51' '0.0320'
Set-TextValue 'E51' '  -0.08%  '
